# Updated cryptos list with refreshed prices and 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.393.31"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "3.695.00"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "693.32"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.84"
$ws.Range("E6").Value = "  -5.36%  "
$ws.Range("D7").Value = "3.694.73"
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -4.64%  "
$ws.Range("E10").Value = "  -8.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -4.52%  "
$ws.Range("E13").Value = "  -5.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.40"
$ws.Range("E14").Value = "  -7.27%  "
$ws.Range("D15").Value = "4.317.42"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "3.698.98"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "69.432.34"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -7.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  -7.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "480.59"
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.98"
$ws.Range("E22").Value = "  -6.52%  "
$ws.Range("E23").Value = "  -7.27%  "
$ws.Range("D25").Value = "3.842.56"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("E26").Value = "  -9.31%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  -8.81%  "
$ws.Range("E30").Value = "  -10.74%  "
$ws.Range("E31").Value = "  -9.81%  "
$ws.Range("E32").Value = "  -7.97%  "
$ws.Range("E33").Value = "  -7.89%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.00"
$ws.Range("E35").Value = "  -7.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "3.663.52"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.49"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.36"
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0931"
$ws.Range("E41").Value = "  -7.89%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.95"
$ws.Range("E45").Value = "  -5.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.94"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.13"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  -15.16%  "

# Rows 49-51 reshuffled: SuiNetwork moves to 51, ONDO moves to 49, FLOKI moves to 50
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.35"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000286"
$ws.Range("E50").Value = "  -8.05%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  -1.31%  "

